$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, shifting existing rows 22-77 down to 23-78.
$ws.Rows.Item(22).Insert()

# Populate the new row 22 with the new price-report record.
$ws.Range("A22").Value = 4
$ws.Range("B22").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C22").Value = "Los Lagos"
$ws.Range("D22").Value = 45272
$ws.Range("E22").Value = 10
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100101
$ws.Range("H22").Value = "Berries"
$ws.Range("I22").Value = 100101001
$ws.Range("J22").Value = "Arándano (blue)"
$ws.Range("K22").Value = "Sin especificar"
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 250
$ws.Range("N22").Value = 5000
$ws.Range("O22").Value = 5000
$ws.Range("P22").Value = 5000
$ws.Range("Q22").Value = "`$/bandeja 12 canastillos 125 gramos"
$ws.Range("R22").Value = "Región del Maule"
$ws.Range("S22").Value = 3333
$ws.Range("T22").Value = 1.5
